# Change the unit test for Sample Annotation validation to suit the new
# version of MSTemplate_Creator: remove the now-unused Transition_Name_Annot
# and ISTD_Annot worksheets, leaving only the Sample_Annot sheet.

$wb = $excel.ActiveWorkbook

$null = $wb.Worksheets("Transition_Name_Annot").Delete()
$null = $wb.Worksheets("ISTD_Annot").Delete()

# Make the remaining (now only) sheet the active one.
$null = $wb.Worksheets("Sample_Annot").Activate()
